$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.39
$ws.Range("H2").Value = 11
$ws.Range("I2").Value = 12
$ws.Range("L2").Value = 1.44
$ws.Range("N2").Value = 3.5
$ws.Range("T2").Value = 2.58
$ws.Range("V2").Value = 1.09
$ws.Range("AA2").Value = 700
$ws.Range("AI2").Value = 260
$ws.Range("AJ2").Value = 10.5
$ws.Range("AN2").Value = 8
$ws.Range("AO2").Value = 560
$ws.Range("F3").Value = 2.92
$ws.Range("G3").Value = 2.94
$ws.Range("I3").Value = 2.76
$ws.Range("L3").Value = 1.4
$ws.Range("W3").Value = 1.51
$ws.Range("AC3").Value = 7.2
$ws.Range("AJ3").Value = 46
$ws.Range("N4").Value = 7.2
$ws.Range("P4").Value = 3.05
$ws.Range("R4").Value = 1.83
$ws.Range("S4").Value = 2.18
$ws.Range("X4").Value = 38
$ws.Range("AB4").Value = 13.5
$ws.Range("AG4").Value = 9.800000000000001
$ws.Range("I5").Value = 3.55
$ws.Range("L5").Value = 1.31
$ws.Range("N5").Value = 4.6
$ws.Range("V5").Value = 1.39
$ws.Range("AB5").Value = 11.5
$ws.Range("AC5").Value = 8.199999999999999
$ws.Range("AJ5").Value = 27
$ws.Range("AN5").Value = 13
$ws.Range("L6").Value = 1.46
$ws.Range("S6").Value = 4.1
$ws.Range("T6").Value = 2.04
$ws.Range("AA6").Value = 140
$ws.Range("AH6").Value = 23
$ws.Range("F7").Value = 3.25
$ws.Range("G7").Value = 3.3
$ws.Range("N7").Value = 3.85
$ws.Range("P7").Value = 1.97
$ws.Range("Q7").Value = 2
$ws.Range("S7").Value = 3.55
$ws.Range("W7").Value = 1.43
$ws.Range("G8").Value = 1.93
$ws.Range("O8").Value = 1.37
$ws.Range("P8").Value = 1.86
$ws.Range("R8").Value = 1.32
$ws.Range("AN8").Value = 15
$ws.Range("J9").Value = 3.8
$ws.Range("K9").Value = 3.85
$ws.Range("L9").Value = 1.37
$ws.Range("AE9").Value = 20
$ws.Range("AG9").Value = 15
$ws.Range("AL9").Value = 55
$ws.Range("F10").Value = 2.56
$ws.Range("I10").Value = 3.45
$ws.Range("J10").Value = 3.15
$ws.Range("K10").Value = 3.2
$ws.Range("L10").Value = 1.5
$ws.Range("Q10").Value = 2.4
$ws.Range("T10").Value = 2.02
$ws.Range("AC10").Value = 7
$ws.Range("AF10").Value = 14
$ws.Range("H11").Value = 3.7
$ws.Range("I11").Value = 3.75
$ws.Range("J11").Value = 3.15
$ws.Range("K11").Value = 3.2
$ws.Range("L11").Value = 1.5
$ws.Range("O11").Value = 1.48
$ws.Range("U11").Value = 1.93
$ws.Range("V11").Value = 1.36
$ws.Range("J12").Value = 5.3
$ws.Range("K12").Value = 5.4
$ws.Range("P12").Value = 2.38
$ws.Range("Q12").Value = 1.69
$ws.Range("V12").Value = 3.4
$ws.Range("Y12").Value = 9.199999999999999
$ws.Range("I13").Value = 5.8
$ws.Range("K13").Value = 4.1
$ws.Range("L13").Value = 1.36
$ws.Range("Q13").Value = 1.86
$ws.Range("U13").Value = 2.12
$ws.Range("V13").Value = 1.2
$ws.Range("AC13").Value = 9
$ws.Range("F14").Value = 5.3
$ws.Range("G14").Value = 5.4
$ws.Range("H14").Value = 1.73
$ws.Range("I14").Value = 1.74
$ws.Range("L14").Value = 1.33
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 4.4
$ws.Range("V14").Value = 2.34
$ws.Range("X14").Value = 18
$ws.Range("AC14").Value = 8.800000000000001
